$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete trailing rows 18-20 (trials 17-19) that are no longer part of the order
$ws.Rows("18:20").Delete()

# Add the new "ITI" header in column D
$ws.Range("D1").Value = "ITI"

# Update ConditionType (column C) and add ITI (column D) values for each trial row
$data = @(
    @{ Row = 2;  C = 4; D = 6 },
    @{ Row = 3;  C = 4; D = 6 },
    @{ Row = 4;  C = 2; D = 7 },
    @{ Row = 5;  C = 2; D = 9 },
    @{ Row = 6;  C = 4; D = 9 },
    @{ Row = 7;  C = 3; D = 8 },
    @{ Row = 8;  C = 2; D = 6 },
    @{ Row = 9;  C = 3; D = 8 },
    @{ Row = 10; C = 1; D = 8 },
    @{ Row = 11; C = 4; D = 6 },
    @{ Row = 12; C = 2; D = 8 },
    @{ Row = 13; C = 1; D = 7 },
    @{ Row = 14; C = 3; D = 6 },
    @{ Row = 15; C = 3; D = 6 },
    @{ Row = 16; C = 1; D = 6 },
    @{ Row = 17; C = 1; D = 6 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
}

# Leave the selection on the last cell touched (D17), matching the saved state
$ws.Range("D17").Select() | Out-Null
